$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "RandomForest Multi-Max R"
$ws.Range("B2").Value = 0.9548387096774194
$ws.Range("C2").Value = 0.8064516129032258

$ws.Range("A3").Value = "XGBoost Multi-Max R"
$ws.Range("B3").Value = 0.967741935483871
$ws.Range("C3").Value = 0.8419354838709677

$ws.Range("A4").Value = "Logistic Regression Multi-Max R"
$ws.Range("B4").Value = 0.9709677419354839
$ws.Range("C4").Value = 0.7645161290322581

$ws.Range("A5").Value = "Voting Classifier Multi-Max R"
$ws.Range("B5").Value = 0.967741935483871
$ws.Range("C5").Value = 0.8290322580645161

$ws.Range("A6").Value = "Stacking Classifier Multi-Max R"
$ws.Range("B6").Value = 0.964516129032258
$ws.Range("C6").Value = 0.8225806451612904
